# Append 4 new daily rows (252-255) to Sheet1, continuing the existing
# time-series table (columns A-D), matching the formatting of the last
# existing data row (A251), which uses a date style (s="2").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style/format of the last existing row's date cell (A251) down
# into the new date cells A252:A255 so the new rows carry the same
# centered/bold/bordered date formatting as all prior rows.
$ws.Range("A251").Copy($ws.Range("A252:A255"))

# Row 252 - 2021-05-10
$ws.Range("A252").Value2 = 44326
$ws.Range("B252").Value2 = 0
$ws.Range("C252").Value2 = 15
$ws.Range("D252").Value2 = 132.8727079457879

# Row 253 - 2021-05-11
$ws.Range("A253").Value2 = 44327
$ws.Range("B253").Value2 = 2
$ws.Range("C253").Value2 = 16
$ws.Range("D253").Value2 = 141.7308884755071

# Row 254 - 2021-05-12
$ws.Range("A254").Value2 = 44328
$ws.Range("B254").Value2 = 0
$ws.Range("C254").Value2 = 14
$ws.Range("D254").Value2 = 124.0145274160687

# Row 255 - 2021-05-13
$ws.Range("A255").Value2 = 44329
$ws.Range("B255").Value2 = 3
$ws.Range("C255").Value2 = 14
$ws.Range("D255").Value2 = 124.0145274160687
